# Update the "2. Data reporter" block (rows 6-10) with the new contact
# details for the Kyrgyz Republic National Statistical Committee, per the
# commit's re-uploaded workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value  = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value  = "Kalymbetova Yryskan"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# The re-uploaded file also landed with the sheet's active-cell selection
# sitting on B8 instead of B2.
$null = $ws.Range("B8").Select()
